$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the existing "Films/fragments" value from C11 to D11,
# and set C11 to the new "Films;fragments" value.
$ws.Range("D11").Value = $ws.Range("C11").Value2
$ws.Range("C11").Value = "Films;fragments"

# Update the active selection to D12, matching the saved workbook state.
$ws.Range("D12").Select()
